$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (updated classifier metrics)
$ws.Range("B4").Value = 0.6015231667975562
$ws.Range("C4").Value = 0.6220000000000001
$ws.Range("D4").Value = 0.5914127684875644
$ws.Range("E4").Value = 0.5965
$ws.Range("F4").Value = 0.63333111688227
$ws.Range("G4").Value = 0.656
$ws.Range("H4").Value = 0.6165133457027109
$ws.Range("I4").Value = 0.6185
$ws.Range("J4").Value = 0.5111399538790977
$ws.Range("K4").Value = 0.513
$ws.Range("L4").Value = 0.5150948030852934
$ws.Range("M4").Value = 0.5165
$ws.Range("N4").Value = 0.6308840169499776
$ws.Range("O4").Value = 0.6609999999999999
$ws.Range("P4").Value = 0.6119448440155988
$ws.Range("Q4").Value = 0.6144999999999999
$ws.Range("R4").Value = 0.6032757859987367
$ws.Range("S4").Value = 0.619
$ws.Range("T4").Value = 0.5978477645633143
$ws.Range("U4").Value = 0.601
$ws.Range("V4").Value = 0.6270932346735956
$ws.Range("W4").Value = 0.649
$ws.Range("X4").Value = 0.6111703256818759
$ws.Range("Y4").Value = 0.6135
$ws.Range("Z4").Value = 0.632592589492169
$ws.Range("AA4").Value = 0.662
$ws.Range("AB4").Value = 0.6150457195910861
$ws.Range("AC4").Value = 0.6165
# Row 5 (updated classifier metrics)
$ws.Range("B5").Value = 0.6522627909298521
$ws.Range("C5").Value = 0.842
$ws.Range("D5").Value = 0.5340715536244627
$ws.Range("E5").Value = 0.554
$ws.Range("F5").Value = 0.6643208308976678
$ws.Range("G5").Value = 0.8530000000000001
$ws.Range("H5").Value = 0.5499871638906561
$ws.Range("I5").Value = 0.573
$ws.Range("J5").Value = 0.6745290942310661
$ws.Range("K5").Value = 0.969
$ws.Range("L5").Value = 0.5178499292574488
$ws.Range("M5").Value = 0.5325
$ws.Range("N5").Value = 0.6629241030651151
$ws.Range("O5").Value = 0.85
$ws.Range("P5").Value = 0.5466654128797196
$ws.Range("Q5").Value = 0.57
$ws.Range("R5").Value = 0.6621633841833641
$ws.Range("S5").Value = 0.8779999999999999
$ws.Range("T5").Value = 0.5328773928810653
$ws.Range("U5").Value = 0.5535
$ws.Range("V5").Value = 0.6752951959547212
$ws.Range("W5").Value = 0.9229999999999998
$ws.Range("X5").Value = 0.5347058109081386
$ws.Range("Y5").Value = 0.5574999999999999
$ws.Range("Z5").Value = 0.6645414853783831
$ws.Range("AA5").Value = 0.875
$ws.Range("AB5").Value = 0.5378070532789861
$ws.Range("AC5").Value = 0.5599999999999999
# Row 6 (updated classifier metrics)
$ws.Range("B6").Value = 0.5745815448526939
$ws.Range("C6").Value = 0.5660000000000001
$ws.Range("D6").Value = 0.5998115209111613
$ws.Range("E6").Value = 0.598
$ws.Range("F6").Value = 0.674261384613661
$ws.Range("G6").Value = 0.6919999999999999
$ws.Range("H6").Value = 0.6771920063597272
$ws.Range("I6").Value = 0.667
$ws.Range("J6").Value = 0.5308052570546892
$ws.Range("K6").Value = 0.525
$ws.Range("L6").Value = 0.5454438940749439
$ws.Range("M6").Value = 0.542
$ws.Range("N6").Value = 0.6542885348579179
$ws.Range("O6").Value = 0.6450000000000001
$ws.Range("P6").Value = 0.6858325444875175
$ws.Range("Q6").Value = 0.6685
$ws.Range("R6").Value = 0.5886385471989254
$ws.Range("S6").Value = 0.5830000000000001
$ws.Range("T6").Value = 0.6142079792501705
$ws.Range("U6").Value = 0.6114999999999999
$ws.Range("V6").Value = 0.6735301340586225
$ws.Range("W6").Value = 0.6970000000000001
$ws.Range("X6").Value = 0.6759493284307114
$ws.Range("Y6").Value = 0.663
$ws.Range("Z6").Value = 0.647515729959139
$ws.Range("AA6").Value = 0.64
$ws.Range("AB6").Value = 0.6765137304929113
$ws.Range("AC6").Value = 0.6585
